$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Metadata" ---
$ws1 = $wb.Worksheets.Item(1)

# Capture existing values of the rows that need to move down by one
# (to make room for the new "Jurisdiction" row) before overwriting anything.
$a15 = $ws1.Range("A15").Value2
$b15 = $ws1.Range("B15").Value2
$a14 = $ws1.Range("A14").Value2
$b14 = $ws1.Range("B14").Value2
$a13 = $ws1.Range("A13").Value2
$b13 = $ws1.Range("B13").Value2
$a12 = $ws1.Range("A12").Value2
$b12 = $ws1.Range("B12").Value2

# Shift rows 12-15 down to rows 13-16
$ws1.Range("A16").Value = $a15
$ws1.Range("B16").Value = $b15
$ws1.Range("A15").Value = $a14
$ws1.Range("B15").Value = $b14
$ws1.Range("A14").Value = $a13
$ws1.Range("B14").Value = $b13
$ws1.Range("A13").Value = $a12
$ws1.Range("B13").Value = $b12

# New "Jurisdiction" row (row 12)
$ws1.Range("A12").Value = "Jurisdiction"
$ws1.Range("B12").Value = ""

# Update the other changed metadata values
$ws1.Range("B3").Value = "0.1.7"
$ws1.Range("B6").Value = "draft"
$ws1.Range("B8").Value = "2024-08-23T10:17:11-05:00"
$ws1.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws1.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"
